$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F3").Value = 45
$sheet1.Range("F6").Value = 1071
$sheet1.Range("F7").Value = 1044
$sheet1.Range("F8").Value = 8120
$sheet1.Range("F9").Value = 135
$sheet1.Range("F10").Value = 204
$sheet1.Range("F11").Value = 6872
$sheet1.Range("F12").Value = 166
$sheet1.Range("F14").Value = 4957
$sheet1.Range("F16").Value = 5380
$sheet1.Range("F17").Value = 1071
$sheet1.Range("F18").Value = 327
$sheet1.Range("F20").Value = 457
$sheet1.Range("F23").Value = 65
$sheet1.Range("F26").Value = 9115
$sheet1.Range("F28").Value = 1645
$sheet1.Range("F29").Value = 645
$sheet1.Range("F32").Value = 1423
$sheet1.Range("F34").Value = 77
$sheet1.Range("F36").Value = 1007
$sheet1.Range("F37").Value = 1863
$sheet1.Range("F38").Value = 240
$sheet1.Range("F40").Value = 4753
$sheet1.Range("F42").Value = 1160
$sheet1.Range("F43").Value = 73
$sheet1.Range("F44").Value = 147
$sheet1.Range("F45").Value = 73
$sheet1.Range("F46").Value = 35
$sheet1.Range("F49").Value = 39
$sheet1.Range("F50").Value = 62

$sheet2 = $wb.Worksheets.Item("演出")
$sheet2.Range("F9").Value = 182

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F4").Value = 45
$sheet4.Range("F8").Value = 1071
$sheet4.Range("F9").Value = 1044
$sheet4.Range("F10").Value = 8120
$sheet4.Range("F11").Value = 135
$sheet4.Range("F12").Value = 204
$sheet4.Range("F13").Value = 6872
$sheet4.Range("F14").Value = 166
$sheet4.Range("F17").Value = 4957
$sheet4.Range("F19").Value = 5380
$sheet4.Range("F20").Value = 1071
$sheet4.Range("F21").Value = 327
$sheet4.Range("F23").Value = 457
$sheet4.Range("F26").Value = 182
$sheet4.Range("F27").Value = 9115
$sheet4.Range("F29").Value = 1645
$sheet4.Range("F30").Value = 645
$sheet4.Range("F33").Value = 1428
$sheet4.Range("F35").Value = 77
$sheet4.Range("F37").Value = 1007
$sheet4.Range("F38").Value = 1863
$sheet4.Range("F39").Value = 240
$sheet4.Range("F41").Value = 4753
$sheet4.Range("F43").Value = 1160
$sheet4.Range("F44").Value = 73
$sheet4.Range("F45").Value = 147
$sheet4.Range("F46").Value = 73
$sheet4.Range("F49").Value = 39
$sheet4.Range("F50").Value = 62
